# Appends 4 new species-observation rows (9-12) to the "Artfynd" sheet,
# extending the used range from A1:AY8 to A1:AY12.
#
# Notes on technique:
#  - Cells that must stay literal text (not be auto-parsed as a date/number
#    by Excel, e.g. "2023-09-24") are written with a leading "'" (force-text
#    prefix), then immediately re-styled to "Normal" so the apostrophe is
#    consumed and no stray quote-prefix formatting lingers on the cell.
#  - Columns that are empty-but-present in the source (I, K, AT, AY) are
#    written the same way with just "'" -> an empty literal text value,
#    matching the original file's empty inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Cells.Item(9, 1).Value = 112292298
$ws.Cells.Item(9, 2).Value = 90843
$ws.Cells.Item(9, 3).Value = 'Ovaliderad'
$ws.Cells.Item(9, 4).Value = 'NT'
$ws.Cells.Item(9, 5).Value = 5448
$ws.Cells.Item(9, 6).Value = 'Svartvit taggsvamp'
$ws.Cells.Item(9, 7).Value = 'Phellodon connatus'
$ws.Cells.Item(9, 8).Value = '(Schultz) nom.prov'
$ws.Cells.Item(9, 9).Value = "'"
$ws.Cells.Item(9, 9).Style = "Normal"
$ws.Cells.Item(9, 11).Value = "'"
$ws.Cells.Item(9, 11).Style = "Normal"
$ws.Cells.Item(9, 16).Value = 'Umeå (Umeå), Vb'
$ws.Cells.Item(9, 17).Value = 763579
$ws.Cells.Item(9, 18).Value = 7089646
$ws.Cells.Item(9, 19).Value = 10
$ws.Cells.Item(9, 20).Value = 'Västerbotten'
$ws.Cells.Item(9, 21).Value = 'Umeå'
$ws.Cells.Item(9, 22).Value = 'Västerbotten'
$ws.Cells.Item(9, 23).Value = 'Umeå socken'
$ws.Cells.Item(9, 25).Value = "'2023-09-24"
$ws.Cells.Item(9, 25).Style = "Normal"
$ws.Cells.Item(9, 26).Value = '16:25'
$ws.Cells.Item(9, 27).Value = "'2023-09-24"
$ws.Cells.Item(9, 27).Style = "Normal"
$ws.Cells.Item(9, 28).Value = '16:25'
$ws.Cells.Item(9, 30).Value = $False
$ws.Cells.Item(9, 31).Value = $False
$ws.Cells.Item(9, 33).Value = $False
$ws.Cells.Item(9, 46).Value = "'"
$ws.Cells.Item(9, 46).Style = "Normal"
$ws.Cells.Item(9, 49).Value = 'André Larencranz'
$ws.Cells.Item(9, 50).Value = 'André Larencranz'
$ws.Cells.Item(9, 51).Value = "'"
$ws.Cells.Item(9, 51).Style = "Normal"

# Row 10
$ws.Cells.Item(10, 1).Value = 112344251
$ws.Cells.Item(10, 2).Value = 90843
$ws.Cells.Item(10, 3).Value = 'Ovaliderad'
$ws.Cells.Item(10, 4).Value = 'NT'
$ws.Cells.Item(10, 5).Value = 5448
$ws.Cells.Item(10, 6).Value = 'Svartvit taggsvamp'
$ws.Cells.Item(10, 7).Value = 'Phellodon connatus'
$ws.Cells.Item(10, 8).Value = '(Schultz) nom.prov'
$ws.Cells.Item(10, 9).Value = "'"
$ws.Cells.Item(10, 9).Style = "Normal"
$ws.Cells.Item(10, 11).Value = "'"
$ws.Cells.Item(10, 11).Style = "Normal"
$ws.Cells.Item(10, 16).Value = 'Umeå (Umeå), Vb'
$ws.Cells.Item(10, 17).Value = 763530
$ws.Cells.Item(10, 18).Value = 7089425
$ws.Cells.Item(10, 19).Value = 10
$ws.Cells.Item(10, 20).Value = 'Västerbotten'
$ws.Cells.Item(10, 21).Value = 'Umeå'
$ws.Cells.Item(10, 22).Value = 'Västerbotten'
$ws.Cells.Item(10, 23).Value = 'Umeå socken'
$ws.Cells.Item(10, 25).Value = "'2023-09-27"
$ws.Cells.Item(10, 25).Style = "Normal"
$ws.Cells.Item(10, 26).Value = '12:29'
$ws.Cells.Item(10, 27).Value = "'2023-09-27"
$ws.Cells.Item(10, 27).Style = "Normal"
$ws.Cells.Item(10, 28).Value = '12:29'
$ws.Cells.Item(10, 30).Value = $False
$ws.Cells.Item(10, 31).Value = $False
$ws.Cells.Item(10, 33).Value = $False
$ws.Cells.Item(10, 46).Value = "'"
$ws.Cells.Item(10, 46).Style = "Normal"
$ws.Cells.Item(10, 49).Value = 'André Larencranz'
$ws.Cells.Item(10, 50).Value = 'André Larencranz'
$ws.Cells.Item(10, 51).Value = "'"
$ws.Cells.Item(10, 51).Style = "Normal"

# Row 11
$ws.Cells.Item(11, 1).Value = 112344188
$ws.Cells.Item(11, 2).Value = 90800
$ws.Cells.Item(11, 3).Value = 'Ovaliderad'
$ws.Cells.Item(11, 4).Value = 'LC'
$ws.Cells.Item(11, 5).Value = 4364
$ws.Cells.Item(11, 6).Value = 'Dropptaggsvamp'
$ws.Cells.Item(11, 7).Value = 'Hydnellum ferrugineum'
$ws.Cells.Item(11, 8).Value = '(Fr.:Fr.) P. Karst.'
$ws.Cells.Item(11, 9).Value = "'"
$ws.Cells.Item(11, 9).Style = "Normal"
$ws.Cells.Item(11, 11).Value = "'"
$ws.Cells.Item(11, 11).Style = "Normal"
$ws.Cells.Item(11, 16).Value = 'Umeå (Umeå), Vb'
$ws.Cells.Item(11, 17).Value = 763543
$ws.Cells.Item(11, 18).Value = 7089426
$ws.Cells.Item(11, 19).Value = 10
$ws.Cells.Item(11, 20).Value = 'Västerbotten'
$ws.Cells.Item(11, 21).Value = 'Umeå'
$ws.Cells.Item(11, 22).Value = 'Västerbotten'
$ws.Cells.Item(11, 23).Value = 'Umeå socken'
$ws.Cells.Item(11, 25).Value = "'2023-09-27"
$ws.Cells.Item(11, 25).Style = "Normal"
$ws.Cells.Item(11, 26).Value = '12:25'
$ws.Cells.Item(11, 27).Value = "'2023-09-27"
$ws.Cells.Item(11, 27).Style = "Normal"
$ws.Cells.Item(11, 28).Value = '12:25'
$ws.Cells.Item(11, 30).Value = $False
$ws.Cells.Item(11, 31).Value = $False
$ws.Cells.Item(11, 33).Value = $False
$ws.Cells.Item(11, 46).Value = "'"
$ws.Cells.Item(11, 46).Style = "Normal"
$ws.Cells.Item(11, 49).Value = 'André Larencranz'
$ws.Cells.Item(11, 50).Value = 'André Larencranz'
$ws.Cells.Item(11, 51).Value = "'"
$ws.Cells.Item(11, 51).Style = "Normal"

# Row 12
$ws.Cells.Item(12, 1).Value = 112344211
$ws.Cells.Item(12, 2).Value = 90823
$ws.Cells.Item(12, 3).Value = 'Ovaliderad'
$ws.Cells.Item(12, 4).Value = 'NT'
$ws.Cells.Item(12, 5).Value = 5966
$ws.Cells.Item(12, 6).Value = 'Motaggsvamp'
$ws.Cells.Item(12, 7).Value = 'Sarcodon squamosus'
$ws.Cells.Item(12, 8).Value = '(Schaeff.) Quél.'
$ws.Cells.Item(12, 9).Value = "'"
$ws.Cells.Item(12, 9).Style = "Normal"
$ws.Cells.Item(12, 11).Value = "'"
$ws.Cells.Item(12, 11).Style = "Normal"
$ws.Cells.Item(12, 16).Value = 'Umeå (Umeå), Vb'
$ws.Cells.Item(12, 17).Value = 763527
$ws.Cells.Item(12, 18).Value = 7089456
$ws.Cells.Item(12, 19).Value = 10
$ws.Cells.Item(12, 20).Value = 'Västerbotten'
$ws.Cells.Item(12, 21).Value = 'Umeå'
$ws.Cells.Item(12, 22).Value = 'Västerbotten'
$ws.Cells.Item(12, 23).Value = 'Umeå socken'
$ws.Cells.Item(12, 25).Value = "'2023-09-27"
$ws.Cells.Item(12, 25).Style = "Normal"
$ws.Cells.Item(12, 26).Value = '12:26'
$ws.Cells.Item(12, 27).Value = "'2023-09-27"
$ws.Cells.Item(12, 27).Style = "Normal"
$ws.Cells.Item(12, 28).Value = '12:26'
$ws.Cells.Item(12, 30).Value = $False
$ws.Cells.Item(12, 31).Value = $False
$ws.Cells.Item(12, 33).Value = $False
$ws.Cells.Item(12, 46).Value = "'"
$ws.Cells.Item(12, 46).Style = "Normal"
$ws.Cells.Item(12, 49).Value = 'André Larencranz'
$ws.Cells.Item(12, 50).Value = 'André Larencranz'
$ws.Cells.Item(12, 51).Value = "'"
$ws.Cells.Item(12, 51).Style = "Normal"

Write-Host "Added rows 9-12 to 'Artfynd' sheet"
